# fix parameterization, run strategy and influence experiments
#
# The "upper" sheet's stakeholder-group rows are reorganized into a new
# run order: the last group of 9 rows (Grower advocacy groups ... environmental
# orgs) moves to the top, the first group of 9 rows (Bureau of Reclamation ...
# Division of Water Quality (SWRCB)) moves to the middle, and the middle
# group of 9 rows (NRCS ... Regional water management groups) moves to the
# bottom. Also update sheet selection / active-tab view state to match the
# reviewed workbook.

$wb    = $excel.ActiveWorkbook
$lower = $wb.Worksheets.Item("lower")
$upper = $wb.Worksheets.Item("upper")

# --- Rearrange the three 9-row blocks on the "upper" sheet ---
# Before: rows 2-10 = Block A, rows 11-19 = Block B, rows 20-28 = Block C
# After : rows 2-10 = Block C, rows 11-19 = Block A, rows 20-28 = Block B
$blockA = $upper.Range("A2:B10").Value2
$blockB = $upper.Range("A11:B19").Value2
$blockC = $upper.Range("A20:B28").Value2

$upper.Range("A2:B10").Value2   = $blockC
$upper.Range("A11:B19").Value2  = $blockA
$upper.Range("A20:B28").Value2  = $blockB

# --- Update view state: "lower" is scrolled and selected but no longer the
#     active tab; "upper" becomes the active tab, scrolled, with a new
#     selection ---
$lower.Activate() | Out-Null
$lower.Range("A11").Select() | Out-Null

$upper.Activate() | Out-Null
$upper.Range("B12").Select() | Out-Null
